$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(167).Insert()

$ws.Range("A167").Value = 10
$ws.Range("B167").Value = "Vega Modelo de Temuco"
$ws.Range("C167").Value = "La Araucanía"
$ws.Range("D167").Value = 45062
$ws.Range("E167").Value = 9
$ws.Range("F167").Value = 100114007
$ws.Range("G167").Value = "Jengibre"
$ws.Range("H167").Value = "Sin especificar"
$ws.Range("I167").Value = "Primera"
$ws.Range("J167").Value = 30
$ws.Range("K167").Value = 24000
$ws.Range("L167").Value = 24000
$ws.Range("M167").Value = 24000
$ws.Range("N167").Value = "$/caja 13 kilos"
$ws.Range("O167").Value = "Perú"
$ws.Range("P167").Value = 1846
$ws.Range("Q167").Value = 13
$ws.Range("R167").Value = "Hortaliza"
